$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(5, "Macroferia Regional de Talca", "Maule", 44890, 7, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Bing",    "Primera", 50, 8000, 8000, 8000, "`$/bandeja 10 kilos", "Provincia de Curicó", 800, 10),
    @(5, "Macroferia Regional de Talca", "Maule", 44890, 7, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Lapins",  "Primera", 50, 8000, 8000, 8000, "`$/bandeja 10 kilos", "Provincia de Curicó", 800, 10),
    @(5, "Macroferia Regional de Talca", "Maule", 44890, 7, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Lapins",  "Segunda", 40, 7000, 7000, 7000, "`$/bandeja 10 kilos", "Provincia de Curicó", 700, 10),
    @(5, "Macroferia Regional de Talca", "Maule", 44890, 7, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Santina", "Primera", 40, 8000, 8000, 8000, "`$/bandeja 10 kilos", "Provincia de Curicó", 800, 10),
    @(5, "Macroferia Regional de Talca", "Maule", 44890, 7, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Santina", "Segunda", 30, 6000, 6000, 6000, "`$/bandeja 10 kilos", "Provincia de Curicó", 600, 10)
)

$startRow = 213
for ($r = 0; $r -lt $newRows.Length; $r++) {
    $rowValues = $newRows[$r]
    $targetRow = $startRow + $r
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($targetRow, $c + 1).Value = $rowValues[$c]
    }
    # Match the date formatting style used by the existing date column (D)
    $ws.Cells.Item($targetRow, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
